$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.9753132958449271
$ws.Range("D2").Value = 0.03331992211473533
$ws.Range("E2").Value = 0.3424319813764996
$ws.Range("F2").Value = 0.8193982407045368
$ws.Range("G2").Value = 0.002439946514007917
$ws.Range("K2").Value = 0.4212605303673342
$ws.Range("L2").Value = 0.1333941180518678
$ws.Range("O2").Value = 2.857116685996772

# Row 3
$ws.Range("B3").Value = 0.9491538420964787
$ws.Range("D3").Value = 0.0316309631125975
$ws.Range("E3").Value = 0.3456736679843964
$ws.Range("F3").Value = 0.8172262884772081
$ws.Range("G3").Value = 0.002442555695383798
$ws.Range("K3").Value = 0.3704118873006621
$ws.Range("L3").Value = 0.1226253742858887
$ws.Range("O3").Value = 2.864074215073288

# Row 4
$ws.Range("B4").Value = 0.9335441425416207
$ws.Range("D4").Value = 0.03058441056181493
$ws.Range("E4").Value = 0.3478060213873846
$ws.Range("F4").Value = 0.8164108716804606
$ws.Range("G4").Value = 0.002444244022874557
$ws.Range("K4").Value = 0.3390339527621506
$ws.Range("L4").Value = 0.116058604185767
$ws.Range("O4").Value = 2.870070639630171

# Row 5
$ws.Range("B5").Value = 0.9272972728535933
$ws.Range("D5").Value = 0.03015556080625004
$ws.Range("E5").Value = 0.3487106935762103
$ws.Range("F5").Value = 0.8162088370362426
$ws.Range("G5").Value = 0.002444953790356842
$ws.Range("K5").Value = 0.3262084976456379
$ws.Range("L5").Value = 0.1133940718436008
$ws.Range("O5").Value = 2.872947636668414

# Row 6
$ws.Range("B6").Value = 0.9262668968257231
$ws.Range("D6").Value = 0.03008420813225854
$ws.Range("E6").Value = 0.348863071970797
$ws.Range("F6").Value = 0.8161831550903855
$ws.Range("G6").Value = 0.0024450729628576
$ws.Range("K6").Value = 0.324076524456558
$ws.Range("L6").Value = 0.1129523243910882
$ws.Range("O6").Value = 2.87345152978375

# Row 7
$ws.Range("B7").Value = 0.9334594320140468
$ws.Range("D7").Value = 0.03057863651228132
$ws.Range("E7").Value = 0.3478180774433497
$ws.Range("F7").Value = 0.8164076196358252
$ws.Range("G7").Value = 0.002444253506908039
$ws.Range("K7").Value = 0.3388611399170429
$ws.Range("L7").Value = 0.1160226227649588
$ws.Range("O7").Value = 2.87010768528512

# Row 8
$ws.Range("B8").Value = 0.9661999116491415
$ws.Range("D8").Value = 0.03273955983978283
$ws.Range("E8").Value = 0.3435202711259695
$ws.Range("F8").Value = 0.8185417759568097
$ws.Range("G8").Value = 0.002440828291598722
$ws.Range("K8").Value = 0.4037608382500082
$ws.Range("L8").Value = 0.1296717092103279
$ws.Range("O8").Value = 2.859157584940334

# Row 9
$ws.Range("B9").Value = 1.033976342087897
$ws.Range("D9").Value = 0.0369006808469976
$ws.Range("E9").Value = 0.3362172770820973
$ws.Range("F9").Value = 0.8268417317479688
$ws.Range("G9").Value = 0.002434793086831311
$ws.Range("K9").Value = 0.5297610944564042
$ws.Range("L9").Value = 0.1567941721438189
$ws.Range("O9").Value = 2.85138073124628

# Row 10
$ws.Range("B10").Value = 1.085933833380437
$ws.Range("D10").Value = 0.03991036568707074
$ws.Range("E10").Value = 0.331535789340105
$ws.Range("F10").Value = 0.8354550928174831
$ws.Range("G10").Value = 0.002430770420744921
$ws.Range("K10").Value = 0.6215360972941255
$ws.Range("L10").Value = 0.1769370486204025
$ws.Range("O10").Value = 2.854039583176956

# Row 11
$ws.Range("B11").Value = 1.110037216949991
$ws.Range("D11").Value = 0.04126906689613463
$ws.Range("E11").Value = 0.3295542063013048
$ws.Range("F11").Value = 0.8399213124744023
$ws.Range("G11").Value = 0.002429028856313178
$ws.Range("K11").Value = 0.6631091947072605
$ws.Range("L11").Value = 0.1861473032891041
$ws.Range("O11").Value = 2.857072365643887

# Row 12
$ws.Range("B12").Value = 1.119231409181225
$ws.Range("D12").Value = 0.04178205294545023
$ws.Range("E12").Value = 0.3288250947045821
$ws.Range("F12").Value = 0.8416914340227635
$ws.Range("G12").Value = 0.002428382011375872
$ws.Range("K12").Value = 0.6788259852859824
$ws.Range("L12").Value = 0.1896417061123401
$ws.Range("O12").Value = 2.858483321313827

# Row 13
$ws.Range("B13").Value = 1.11724831374238
$ws.Range("D13").Value = 0.04167164044604021
$ws.Range("E13").Value = 0.3289811761957413
$ws.Range("F13").Value = 0.8413066982845834
$ws.Range("G13").Value = 0.002428520759394584
$ws.Range("K13").Value = 0.6754422633338208
$ws.Range("L13").Value = 0.1888888285036927
$ws.Range("O13").Value = 2.858167766573985

# Row 14
$ws.Range("B14").Value = 1.110792293401914
$ws.Range("D14").Value = 0.04131130127240823
$ws.Range("E14").Value = 0.3294937956561128
$ws.Range("F14").Value = 0.8400653605775403
$ws.Range("G14").Value = 0.002428975386670359
$ws.Range("K14").Value = 0.6644027507123837
$ws.Range("L14").Value = 0.186434656634475
$ws.Range("O14").Value = 2.857183182943004

# Row 15
$ws.Range("B15").Value = 1.106846474798516
$ws.Range("D15").Value = 0.04109038379708352
$ws.Range("E15").Value = 0.3298105594309604
$ws.Range("F15").Value = 0.839315277341683
$ws.Range("G15").Value = 0.002429255504575506
$ws.Range("K15").Value = 0.6576373170060776
$ws.Range("L15").Value = 0.1849322723120537
$ws.Range("O15").Value = 2.856614292742307

# Row 16
$ws.Range("B16").Value = 1.084367989611167
$ws.Range("D16").Value = 0.03982135944992393
$ws.Range("E16").Value = 0.3316682680053589
$ws.Range("F16").Value = 0.8351742468076537
$ws.Range("G16").Value = 0.002430886009251765
$ws.Range("K16").Value = 0.6188155949441239
$ws.Range("L16").Value = 0.1763360770413414
$ws.Range("O16").Value = 2.853878100136114

# Row 17
$ws.Range("B17").Value = 1.070697600564472
$ws.Range("D17").Value = 0.03904016667171106
$ws.Range("E17").Value = 0.3328458178741673
$ws.Range("F17").Value = 0.832774251032518
$ws.Range("G17").Value = 0.002431908861740376
$ws.Range("K17").Value = 0.594954155605933
$ws.Range("L17").Value = 0.1710745954977568
$ws.Range("O17").Value = 2.852666746624351

# Row 18
$ws.Range("B18").Value = 1.062878801657945
$ws.Range("D18").Value = 0.03858986607522752
$ws.Range("E18").Value = 0.3335370496012295
$ws.Range("F18").Value = 0.8314454095836936
$ws.Range("G18").Value = 0.002432505500210827
$ws.Range("K18").Value = 0.5812131835484422
$ws.Range("L18").Value = 0.1680527754094641
$ws.Range("O18").Value = 2.85214160253841

# Row 19
$ws.Range("B19").Value = 1.060239070658866
$ws.Range("D19").Value = 0.0384372346606483
$ws.Range("E19").Value = 0.33377348290316
$ws.Range("F19").Value = 0.8310043423375362
$ws.Range("G19").Value = 0.002432708942820746
$ws.Range("K19").Value = 0.5765579152296425
$ws.Range("L19").Value = 0.1670304055874681
$ws.Range("O19").Value = 2.851993259342123

# Row 20
$ws.Range("B20").Value = 1.072148281602693
$ws.Range("D20").Value = 0.03912342750762576
$ws.Range("E20").Value = 0.3327190234709683
$ws.Range("F20").Value = 0.8330243966251487
$ws.Range("G20").Value = 0.002431799116665187
$ws.Range("K20").Value = 0.5974959600891054
$ws.Range("L20").Value = 0.1716342298971938
$ws.Range("O20").Value = 2.852777935527484

# Row 21
$ws.Range("B21").Value = 1.112686774976424
$ws.Range("D21").Value = 0.04141718320908438
$ws.Range("E21").Value = 0.3293426497547038
$ws.Range("F21").Value = 0.8404278307726685
$ws.Range("G21").Value = 0.00242884150894231
$ws.Range("K21").Value = 0.6676460358796135
$ws.Range("L21").Value = 0.1871553259419727
$ws.Range("O21").Value = 2.857465252094443

# Row 22
$ws.Range("B22").Value = 1.139569951213673
$ws.Range("D22").Value = 0.04290738645815395
$ws.Range("E22").Value = 0.3272599689314095
$ws.Range("F22").Value = 0.8457260949196268
$ws.Range("G22").Value = 0.002426982236052597
$ws.Range("K22").Value = 0.713341041473484
$ws.Range("L22").Value = 0.1973381440912902
$ws.Range("O22").Value = 2.862058922936029

# Row 23
$ws.Range("B23").Value = 1.125186467170209
$ws.Range("D23").Value = 0.04211285982196955
$ws.Range("E23").Value = 0.3283601972540033
$ws.Range("F23").Value = 0.8428562286317316
$ws.Range("G23").Value = 0.00242796784133742
$ws.Range("K23").Value = 0.6889669254081809
$ws.Range("L23").Value = 0.1918998581506628
$ws.Range("O23").Value = 2.859467071101477

# Row 24
$ws.Range("B24").Value = 1.071492302547568
$ws.Range("D24").Value = 0.03908578896443515
$ws.Range("E24").Value = 0.3327763028550557
$ws.Range("F24").Value = 0.8329111471066852
$ws.Range("G24").Value = 0.002431848705671329
$ws.Range("K24").Value = 0.5963468809752328
$ws.Range("L24").Value = 0.1713812095988487
$ws.Range("O24").Value = 2.852727133518329

# Row 25
$ws.Range("B25").Value = 1.015260196511463
$ws.Range("D25").Value = 0.03578326359281192
$ws.Range("E25").Value = 0.3380726792111144
$ws.Range("F25").Value = 0.8241551347096063
$ws.Range("G25").Value = 0.002436353226111658
$ws.Range("K25").Value = 0.4958129687741746
$ws.Range("L25").Value = 0.1494187936663991
$ws.Range("O25").Value = 2.852015919175102
